$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# ---- Row 6 ----
$ws.Range("G6").Value  = 2.25
$ws.Range("H6").Value  = 3.1
$ws.Range("I6").Value  = 3.1
$ws.Range("J6").Value  = 2.77
$ws.Range("K6").Value  = 2.1
$ws.Range("L6").Value  = 3.55
$ws.Range("M6").Value  = 1.01
$ws.Range("N6").Value  = 8.1
$ws.Range("O6").Value  = 1.3
$ws.Range("P6").Value  = 2.95
$ws.Range("Q6").Value  = 1.93
$ws.Range("R6").Value  = 1.78
$ws.Range("X6").Value  = 11.5
$ws.Range("Z6").Value  = 23
$ws.Range("AC6").Value = 9.25
$ws.Range("AD6").Value = 6.1
$ws.Range("AF6").Value = 60
$ws.Range("AG6").Value = 9.25
$ws.Range("AH6").Value = 16
$ws.Range("AI6").Value = 10.75
$ws.Range("AK6").Value = 27
$ws.Range("AL6").Value = 35
$ws.Range("AM6").Value = 450
$ws.Range("AN6").Value = 4.2
$ws.Range("AO6").Value = 11.25
$ws.Range("AP6").Value = 17.5
$ws.Range("AR6").Value = 65
$ws.Range("AS6").Value = 200
$ws.Range("AT6").Value = 2.62
$ws.Range("AU6").Value = 6.5
$ws.Range("AV6").Value = 50
$ws.Range("AW6").Value = 5.1
$ws.Range("AX6").Value = 16.5
$ws.Range("AZ6").Value = 80

# ---- Row 15 ----
$ws.Range("G15").Value  = 2.4
$ws.Range("I15").Value  = 3.25
$ws.Range("J15").Value  = 3.2
$ws.Range("L15").Value  = 4
$ws.Range("M15").Value  = 1.1
$ws.Range("N15").Value  = 7
$ws.Range("X15").Value  = 10
$ws.Range("AI15").Value = 13

# ---- Row 16 ----
$ws.Range("M16").Value = 1.06
$ws.Range("N16").Value = 10
$ws.Range("O16").Value = 1.3
$ws.Range("P16").Value = 3.4
